# Update "want to go" counts (column F) on several sheets to reflect
# the latest data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5399
$wsExpo.Range("F4").Value = 11510
$wsExpo.Range("F5").Value = 283
$wsExpo.Range("F7").Value = 169
$wsExpo.Range("F8").Value = 261
$wsExpo.Range("F9").Value = 999

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 17

# 全部类型 (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 17
$wsAll.Range("F4").Value = 5399
$wsAll.Range("F7").Value = 11510
$wsAll.Range("F8").Value = 283
$wsAll.Range("F10").Value = 169
$wsAll.Range("F13").Value = 261
$wsAll.Range("F14").Value = 999
